$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price (D) and volume (E) columns for rows with changes ---
# Plain text / already-non-numeric-looking values: direct assignment
$ws.Range("D2").Value = "56.438.95"
$ws.Range("E2").Value = "  +3.97%  "
$ws.Range("D3").Value = "2.992.66"
$ws.Range("E3").Value = "  +4.59%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +7.55%  "
$ws.Range("E6").Value = "  +8.96%  "
$ws.Range("E8").Value = "  +7.94%  "
$ws.Range("E9").Value = "  +14.24%  "
$ws.Range("E10").Value = "  +14.02%  "
$ws.Range("E11").Value = "  +7.29%  "
$ws.Range("E12").Value = "  +4.69%  "
$ws.Range("D13").Value = "3.502.79"
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("E14").Value = "  +11.38%  "
$ws.Range("E15").Value = "  +15.46%  "
$ws.Range("D16").Value = "56.482.88"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "2.990.18"
$ws.Range("E17").Value = "  +4.52%  "
$ws.Range("E18").Value = "  +10.28%  "
$ws.Range("E19").Value = "  +9.05%  "
$ws.Range("E20").Value = "  +11.16%  "
$ws.Range("E21").Value = "  +11.10%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +7.63%  "
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  +6.71%  "
$ws.Range("D27").Value = "0.0₃0905"
$ws.Range("E27").Value = "  +13.48%  "
$ws.Range("E28").Value = "  +4.53%  "
$ws.Range("E29").Value = "  +12.11%  "
$ws.Range("E30").Value = "  +8.76%  "
$ws.Range("E31").Value = "  +9.42%  "
$ws.Range("E32").Value = "  +9.11%  "
$ws.Range("E33").Value = "  +16.03%  "
$ws.Range("E34").Value = "  +6.58%  "
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("E36").Value = "  +4.13%  "
$ws.Range("E37").Value = "  +10.03%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "3.024.68"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E40").Value = "  +4.11%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +7.35%  "
$ws.Range("D43").Value = "2.253.36"
$ws.Range("E43").Value = "  +10.52%  "
$ws.Range("E44").Value = "  +7.12%  "
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("E46").Value = "  +6.78%  "
$ws.Range("E47").Value = "  +23.15%  "
$ws.Range("E48").Value = "  +12.00%  "
$ws.Range("E49").Value = "  +9.14%  "
$ws.Range("E50").Value = "  +7.17%  "
$ws.Range("E51").Value = "  +11.06%  "

# Numeric-looking price strings: force text type via NumberFormat "@" trick,
# then ClearFormats so no residual style is left on the cell (matches original,
# un-styled inline-string cells) while the stored type/value stay textual.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.57"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.11"
$ws.Range("D6").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.53"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.352"
$ws.Range("D11").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.64"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000154"
$ws.Range("D15").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.82"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.43"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.82"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.33"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.474"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.34"
$ws.Range("D24").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.48"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.90"
$ws.Range("D29").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.52"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.05"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.58"
$ws.Range("D35").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0673"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.33"
$ws.Range("D38").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.43"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.644"
$ws.Range("D42").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.988"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.59"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.97"
$ws.Range("D47").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.98"
$ws.Range("D50").ClearFormats()
